# "Generate Report for Archive"
#
# The localization status report is regenerated:
#   * The status text for the sample file changes from
#     "Ready for handoff" to "In Translation" on every sheet that
#     references it (Overview!E2/F2, zh-cn!C2, de-de!C2).
#   * The Status column on the Overview, zh-cn and de-de sheets is
#     narrowed to fit the new (shorter) status text.
#
# Note: Excel's ColumnWidth property only stores values that are
# quantized to whole-pixel increments (i.e. multiples of 1/6 of a
# character at the default font), so the COM-settable width that ends
# up closest to the recorded OOXML width (13.4101845877511 characters)
# is used below.

$wb = $excel.ActiveWorkbook

$targetColumnWidth = 12.5   # rounds to the stored width closest to 13.4101845877511
$newStatus = "In Translation"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $targetColumnWidth   # column E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = $targetColumnWidth   # column F (de-de)

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $targetColumnWidth       # column C (Status)

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $targetColumnWidth       # column C (Status)
